$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the Rx / Ma tuy entry) - this shifts row 3 up to become the new row 2
$ws.Rows("2:2").Delete()
